# Add new "absolute order" test case row (row 17) to Sheet1:
# F17=S_1, G17=V_2, H17=O1_3, I17=O2_4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F17").Value = "S_1"
$ws.Range("G17").Value = "V_2"
$ws.Range("H17").Value = "O1_3"
$ws.Range("I17").Value = "O2_4"

# Matches the author's final selection left in the saved file (I18, just
# below the newly added row).
[void]$ws.Range("I18").Select()
